$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("D7").Value = 44242
$ws.Range("J7").Value = 60

# Row 8
$ws.Range("D8").Value = 44242
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 10000
$ws.Range("L8").Value = 10000
$ws.Range("M8").Value = 10000
$ws.Range("P8").Value = 556

# Row 9
$ws.Range("D9").Value = 44238
$ws.Range("J9").Value = 90

# Row 10
$ws.Range("D10").Value = 44238
$ws.Range("J10").Value = 80
$ws.Range("K10").Value = 11000
$ws.Range("L10").Value = 11000
$ws.Range("M10").Value = 11000
$ws.Range("P10").Value = 611
